$wb = $excel.ActiveWorkbook

# --- Tabelle1: move the selection (was A95:XFD95) to I19 ---
$tabelle1 = $wb.Worksheets.Item("Tabelle1")
$tabelle1.Range("I19").Select()

# --- add the new "Tabelle3" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Tabelle3"

# A2: upper bound used by the scaling formulas below
$newSheet.Cells.Item(2, 1).Value = 350

# Row 3 (B3:O3): exponents 1..14
for ($c = 2; $c -le 15; $c++) {
    $newSheet.Cells.Item(3, $c).Value = $c - 1
}

# Rows 4..24: column A is the fractional "progress" value, 0..1 in steps of 0.05
# columns B..O hold MAX(100, MIN($A$2, ROUND(100 + (A^exp) * ($A$2-100), 0)))
$Avals = @(0, 0.05, 0.1, 0.15, 0.2, 0.25, 0.3, 0.35, 0.4, 0.45, 0.5, 0.55, 0.6, 0.65, 0.7, 0.75, 0.8, 0.85, 0.9, 0.95, 1)
for ($i = 0; $i -lt $Avals.Count; $i++) {
    $row = 4 + $i
    $newSheet.Cells.Item($row, 1).Value = $Avals[$i]
    for ($c = 2; $c -le 15; $c++) {
        $colLetter = [char](64 + $c)
        $formula = "=MAX(100, MIN(`$A`$2, ROUND(100 + (`$A$row ^ $colLetter`$3) * (`$A`$2 - 100), 0)))"
        $newSheet.Cells.Item($row, $c).Formula = $formula
    }
}

# Tabelle3 ends up the active / selected tab, with A3 selected
$newSheet.Activate()
$newSheet.Range("A3").Select()

Write-Output "done"
